$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": a new daily column ("02-nov") is inserted right
#     before the "01-oct." column (column DG, the 111th column). Excel
#     shifts every following column (formerly DG..EK) one position to the
#     right (DH..EL), preserving their formatting. ---
$ws1 = $wb.Worksheets.Item("Prix Spot")
$ws1.Columns.Item(111).Insert()

$ws1.Range("DG1").Value = "02-nov"
for ($r = 2; $r -le 25; $r++) {
    $ws1.Cells.Item($r, 111).Value = "-"
}

# --- Sheet "Gaz": append the next day's quote as a new row. Force the
#     date column to stay plain text (like the existing rows) instead of
#     letting Excel auto-convert the "yyyy-mm-dd" string into a date. ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A139").NumberFormat = "@"
$ws2.Range("A139").Value = "2025-10-31"
$ws2.Range("A139").ClearFormats()
$ws2.Range("B139").Value = 29.9

# --- Sheet "CO2": same kind of new row appended. ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A139").NumberFormat = "@"
$ws3.Range("A139").Value = "2025-10-31"
$ws3.Range("A139").ClearFormats()
$ws3.Range("B139").Value = 78
